$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they remain text (matches source formatting)
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D13", "D15", "D17", "D20", "D21", "D22", "D24", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D39", "D40", "D42", "D43", "D44", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '42.394.42'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.185.31'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '253.91'
$ws.Range('E5').Value = '  +5.47%  '
$ws.Range('D6').Value = '0.607'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('D7').Value = '74.28'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  -2.25%  '
$ws.Range('D10').Value = '40.66'
$ws.Range('E10').Value = '  -1.54%  '
$ws.Range('D11').Value = '0.0915'
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').Value = '6.79'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('D14').Value = '2.510.31'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = '14.21'
$ws.Range('E15').Value = '  -3.15%  '
$ws.Range('D16').Value = '2.176.86'
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').Value = '0.770'
$ws.Range('E17').Value = '  -3.59%  '
$ws.Range('D18').Value = '42.334.59'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('D20').Value = '70.60'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').Value = '5.88'
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').Value = '227.00'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').Value = '9.42'
$ws.Range('E24').Value = '  -6.78%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('D26').Value = '10.46'
$ws.Range('E26').Value = '  -4.27%  '
$ws.Range('D27').Value = '3.35'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').Value = '170.18'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').Value = '36.60'
$ws.Range('E31').Value = '  +8.89%  '
$ws.Range('D32').Value = '20.03'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('D33').Value = '0.0807'
$ws.Range('E33').Value = '  +1.90%  '
$ws.Range('D34').Value = '5.13'
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('E35').Value = '  -0.81%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').Value = '4.25'
$ws.Range('E37').Value = '  -3.66%  '
$ws.Range('E38').Value = '  +5.01%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '11.84'
$ws.Range('E39').Value = '  -5.83%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '2.06'
$ws.Range('E40').Value = '  -2.97%  '
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('D42').Value = '59.36'
$ws.Range('E42').Value = '  -2.00%  '
$ws.Range('D43').Value = '5.14'
$ws.Range('E43').Value = '  -6.17%  '
$ws.Range('D44').Value = '102.72'
$ws.Range('E44').Value = '  +3.35%  '
$ws.Range('E45').Value = '  +10.77%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '2.44'
$ws.Range('E46').Value = '  +6.60%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '8.25'
$ws.Range('E47').Value = '  -3.58%  '
$ws.Range('D48').Value = '0.0970'
$ws.Range('E48').Value = '  -0.91%  '
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').Value = '1.13'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('E51').Value = '  +0.42%  '
